# Update market price / profit figures pulled by the scheduled pricing runner.
# Each worksheet (one per crafting job) has columns H-N holding fetched market
# prices and computed Leve profits; this script overwrites the refreshed values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 112
$ws.Range("H112").Value = 24526.256
$ws.Range("I112").Value = 763
$ws.Range("J112").Value = 26308.5
$ws.Range("K112").Value = 2289
$ws.Range("L112").Value = 78925.5
$ws.Range("M112").Value = -1181
$ws.Range("N112").Value = -81141.5
# Row 129
$ws.Range("H129").Value = 1076.1389
$ws.Range("J129").Value = 1101.2572
$ws.Range("L129").Value = 3303.7716
$ws.Range("N129").Value = -13303.7716
# Row 137
$ws.Range("H137").Value = 2464.6333
$ws.Range("I137").Value = 2160.8572
$ws.Range("J137").Value = 3817.818
$ws.Range("K137").Value = 6482.571599999999
$ws.Range("L137").Value = 11453.454
$ws.Range("M137").Value = -3932.571599999999
$ws.Range("N137").Value = -16553.454

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 21236
$ws.Range("I61").Value = 26197.5
$ws.Range("J61").Value = 1390
$ws.Range("K61").Value = 26197.5
$ws.Range("L61").Value = 1390
$ws.Range("M61").Value = -25985.5
$ws.Range("N61").Value = -1814
# Row 74
$ws.Range("H74").Value = 8080.2856
$ws.Range("I74").Value = 9092.5
$ws.Range("J74").Value = 2007
$ws.Range("K74").Value = 9092.5
$ws.Range("L74").Value = 2007
$ws.Range("M74").Value = -8218.5
$ws.Range("N74").Value = -3755
# Row 77
$ws.Range("H77").Value = 8080.2856
$ws.Range("I77").Value = 9092.5
$ws.Range("J77").Value = 2007
$ws.Range("K77").Value = 45462.5
$ws.Range("L77").Value = 10035
$ws.Range("M77").Value = -41094.5
$ws.Range("N77").Value = -18771
# Row 88
$ws.Range("H88").Value = 2355.5557
$ws.Range("I88").Value = 2033.3334
$ws.Range("J88").Value = 2516.6667
$ws.Range("K88").Value = 2033.3334
$ws.Range("L88").Value = 2516.6667
$ws.Range("M88").Value = -1627.3334
$ws.Range("N88").Value = -3328.6667
# Row 91
$ws.Range("H91").Value = 2355.5557
$ws.Range("I91").Value = 2033.3334
$ws.Range("J91").Value = 2516.6667
$ws.Range("K91").Value = 2033.3334
$ws.Range("L91").Value = 2516.6667
$ws.Range("M91").Value = -629.3334
$ws.Range("N91").Value = -5324.6667
# Row 132
$ws.Range("H132").Value = 14260.9375
$ws.Range("I132").Value = 9555.166999999999
$ws.Range("J132").Value = 28378.25
$ws.Range("K132").Value = 28665.501
$ws.Range("L132").Value = 85134.75
$ws.Range("M132").Value = -26135.501
$ws.Range("N132").Value = -90194.75
# Row 136
$ws.Range("H136").Value = 21236
$ws.Range("I136").Value = 26197.5
$ws.Range("J136").Value = 1390
$ws.Range("K136").Value = 78592.5
$ws.Range("L136").Value = 4170
$ws.Range("M136").Value = -76042.5
$ws.Range("N136").Value = -9270

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 3295
$ws.Range("I86").Value = 3323.6296
$ws.Range("J86").Value = 3166.1667
$ws.Range("K86").Value = 3323.6296
$ws.Range("L86").Value = 3166.1667
$ws.Range("M86").Value = -2200.6296
$ws.Range("N86").Value = -5412.1667
# Row 89
$ws.Range("H89").Value = 3295
$ws.Range("I89").Value = 3323.6296
$ws.Range("J89").Value = 3166.1667
$ws.Range("K89").Value = 16618.148
$ws.Range("L89").Value = 15830.8335
$ws.Range("M89").Value = -11002.148
$ws.Range("N89").Value = -27062.8335
# Row 134
$ws.Range("H134").Value = 4166.091
$ws.Range("I134").Value = 4359.3667
$ws.Range("J134").Value = 2233.3333
$ws.Range("K134").Value = 13078.1001
$ws.Range("L134").Value = 6699.999899999999
$ws.Range("M134").Value = -10543.1001
$ws.Range("N134").Value = -11769.9999

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 41062.715
$ws.Range("I31").Value = 74199.734
$ws.Range("J31").Value = 2827.6924
$ws.Range("K31").Value = 74199.734
$ws.Range("L31").Value = 2827.6924
$ws.Range("M31").Value = -73904.734
$ws.Range("N31").Value = -3417.6924
# Row 34
$ws.Range("H34").Value = 41062.715
$ws.Range("I34").Value = 74199.734
$ws.Range("J34").Value = 2827.6924
$ws.Range("K34").Value = 74199.734
$ws.Range("L34").Value = 2827.6924
$ws.Range("M34").Value = -73997.734
$ws.Range("N34").Value = -3231.6924
# Row 58
$ws.Range("H58").Value = 1150.4445
$ws.Range("I58").Value = 1294.4762
$ws.Range("J58").Value = 646.3333
$ws.Range("K58").Value = 1294.4762
$ws.Range("L58").Value = 646.3333
$ws.Range("M58").Value = -1091.4762
$ws.Range("N58").Value = -1052.3333
# Row 132
$ws.Range("H132").Value = 9428.440000000001
$ws.Range("I132").Value = 6437.278
$ws.Range("J132").Value = 17120
$ws.Range("K132").Value = 19311.834
$ws.Range("L132").Value = 51360
$ws.Range("M132").Value = -16781.834
$ws.Range("N132").Value = -56420
# Row 134
$ws.Range("H134").Value = 8361.799999999999
$ws.Range("I134").Value = 9374.166999999999
$ws.Range("J134").Value = 7686.8887
$ws.Range("K134").Value = 28122.501
$ws.Range("L134").Value = 23060.6661
$ws.Range("M134").Value = -25587.501
$ws.Range("N134").Value = -28130.6661
# Row 136
$ws.Range("H136").Value = 1150.4445
$ws.Range("I136").Value = 1294.4762
$ws.Range("J136").Value = 646.3333
$ws.Range("K136").Value = 3883.4286
$ws.Range("L136").Value = 1938.9999
$ws.Range("M136").Value = -1333.4286
$ws.Range("N136").Value = -7038.9999

$ws = $wb.Worksheets.Item("CUL")
# Row 56
$ws.Range("H56").Value = 5125
$ws.Range("I56").Value = 5125
$ws.Range("K56").Value = 5125
$ws.Range("M56").Value = -4595

$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 13093.529
$ws.Range("I132").Value = 7661.875
$ws.Range("K132").Value = 22985.625
$ws.Range("M132").Value = -20455.625

$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 4750.171
$ws.Range("I132").Value = 4907.5454
$ws.Range("J132").Value = 4101
$ws.Range("K132").Value = 14722.6362
$ws.Range("L132").Value = 12303
$ws.Range("M132").Value = -12192.6362
$ws.Range("N132").Value = -17363
# Row 136
$ws.Range("H136").Value = 2850.634
$ws.Range("I136").Value = 1964.05
$ws.Range("J136").Value = 3695
$ws.Range("K136").Value = 5892.15
$ws.Range("L136").Value = 11085
$ws.Range("M136").Value = -3342.15
$ws.Range("N136").Value = -16185

$ws = $wb.Worksheets.Item("WVR")
# Row 24
$ws.Range("H24").Value = 70000
$ws.Range("J24").Value = 70000
$ws.Range("L24").Value = 70000
$ws.Range("N24").Value = -70460
# Row 132
$ws.Range("H132").Value = 15027.7
$ws.Range("I132").Value = 19029.834
$ws.Range("K132").Value = 57089.50199999999
$ws.Range("M132").Value = -54559.50199999999
# Row 136
$ws.Range("H136").Value = 37042860
$ws.Range("I136").Value = 47625884
$ws.Range("J136").Value = 2265
$ws.Range("K136").Value = 142877652
$ws.Range("L136").Value = 6795
$ws.Range("M136").Value = -142875102
$ws.Range("N136").Value = -11895

